$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Remove the empty "Title 1" placeholder shape, leaving only "TextBox 2"
foreach ($shp in @($s.Shapes)) {
    if ($shp.Name -eq "Title 1") {
        $shp.Delete()
    }
}
